# Update schedule and 解答 (commit: "update schedule and 解答")
#
# Row 12: fill in the completion date (C12) that was previously blank.
# Row 13: a brand new schedule entry (task/date/notes) that used to be
#         an empty placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- C12: record the completion date for the existing row 12 task ----
# Copy the date formatting from the neighbouring date cell (B12) so the
# new value picks up the same number format / style, then write the value.
$ws.Range("B12").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = 43223

# ---- Row 13: new task entry ----
# A13 - task description (wrap-text style copied from A12)
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "1. Fragment的使用(静态添加(定义在xml中)，动态添加，与Activity间的通信，生命周期)`n2. Linux命令(find, file, touch, ifconfig, clear, export, reboot(重启命令)，sudo. exit)`n"

# B13 - expected date (date style copied from B12)
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = 43225

# D13 - notes / 解答 (wrap-text style copied from D12)
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").Value = "第一行代码第4章，Fragment使用的步骤要记录好。`nShell的基本命令就到这，后面学shell脚本时有碰到新的命令再学"

$excel.CutCopyMode = 0

# Move the active selection to A13 (matches the saved view state).
$ws.Activate() | Out-Null
$ws.Range("A13").Select() | Out-Null
